$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's rule table (B3:E11) records a 4th rule "R40" / 22 / 23 /
# "Good Night" in row 11. This edit rewrites the rule name in B11 to the
# literal text "1" (kept as a text value, not a number), leaving the
# cell's existing style/formatting untouched.
#
# A plain "$ws.Range('B11').Value = '1'" would be auto-recognised as a
# number by Excel's input parser, which is not what we want here - the
# saved value must stay a text string. So the new text is staged on a
# scratch cell well outside the used range (B3:E11), using a formula
# that evaluates to the text "1" (this avoids touching any cell
# NumberFormat/style), then copied across with Paste Special > Values
# so only B11's stored value/type changes while its style (s="23")
# stays exactly as it was.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""1"""
$scratch.Copy()

$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = $false
$scratch.Clear()
